# project 5b — extend the graph-map counting table (P/S columns) down through
# row 44 and add a COUNT() summary row, per the "working on project 5b" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minimize the workbook window (workbookView minimized="1" in the diff).
try { $excel.WindowState = -4140 } catch { }

# New P/S data rows 29-44.
$rows = @(
    @(29, 0,  0),
    @(30, 1,  1),
    @(31, 2,  2),
    @(32, 3,  3),
    @(33, 4,  10),
    @(34, 5,  14),
    @(35, 6,  20),
    @(36, 7,  23),
    @(37, 12, 24),
    @(38, 16, 25),
    @(39, 21, 26),
    @(40, 27, 27),
    @(41, 28, 28),
    @(42, 29, 29),
    @(43, 30, 30),
    @(44, 31, 31)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 16).Value = $row[1]   # column P
    $ws.Cells.Item($r, 19).Value = $row[2]   # column S
}

# Summary row 46: label, COUNT() formula, and the matching literal total.
$ws.Range("N46").Value = "Count:"
$ws.Range("P46").Formula = "=COUNT(P29:P44)"
$ws.Range("S46").Value = 16

# Move the selection to match the author's last-saved cursor position.
[void]$ws.Range("S47").Select()
